$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E13 value (corrected num_matches for season 11)
$ws.Range("E13").Value = 1182353

# Add new season row 14 (season index 12, "M3_01 Wolf 2021")
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "M3_01 Wolf 2021"
$ws.Range("C14").Value = 9637
$ws.Range("D14").Value = 10653
$ws.Range("E14").Value = 808651
$ws.Range("F14").Value = 9916
$ws.Range("G14").Value = 10044
$ws.Range("H14").Value = 10295

# Add new season row 15 (season index 13, "M3_02 Love 2021")
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "M3_02 Love 2021"
$ws.Range("C15").Value = 9684
$ws.Range("D15").Value = 10714
$ws.Range("E15").Value = 917491
$ws.Range("F15").Value = 9975
$ws.Range("G15").Value = 10097
$ws.Range("H15").Value = 10325
